$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing response in row 2 (column C): was 1, should be 0
$ws.Range("C2").Value = 0.0

# New survey response rows (rows 3 and 4), mirroring the header order in row 1
$row3 = @(1.0,1.0,1.0,0.0,1.0,0.0,1.0,0.0,1.0,0.0,1.0,0.0,1.0,0.0,1.0,2.0,1.0,0.0,1.0,1.0,1.0,2.0,1.0,0.0,1.0,2.0,0.0,1.0,0.0,1.0,1.0,1.0,0.0,1.0)
$row4 = @(1.0,1.0,1.0,0.0,1.0,1.0,1.0,0.0,1.0,0.0,1.0,0.0,1.0,0.0,1.0,2.0,1.0,0.0,1.0,1.0,1.0,2.0,1.0,0.0,1.0,2.0,0.0,1.0,0.0,1.0,1.0,1.0,0.0,1.0)

for ($col = 1; $col -le 34; $col++) {
    $ws.Cells.Item(3, $col).Value = $row3[$col - 1]
    $ws.Cells.Item(4, $col).Value = $row4[$col - 1]
}
